# Adicao de texto com dados da turma
#
# Slide 1, shape "CaixaDeTexto 16" (the "DISCIPLINA:" box): the second
# paragraph ("QUALIDADE DE SOFTWARE E GOVERNANCA DE TI") gets a new
# paragraph right after it containing "TURMA3SI". The textbox uses
# spAutoFit, so growing the text automatically grows p:spPr/a:xfrm/a:ext.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(7)

# sanity check - make sure we grabbed the right shape
if ($sh.Name -ne "CaixaDeTexto 16") {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq "CaixaDeTexto 16") {
            $sh = $s.Shapes.Item($i)
        }
    }
}

$tr = $sh.TextFrame.TextRange

# Paragraph 2 is "QUALIDADE DE SOFTWARE E GOVERNANCA DE TI". Insert a new
# paragraph right after it (leading carriage return starts a new
# paragraph) holding the "TURMA3SI" run; the new run inherits the bold
# it-IT run formatting of the paragraph it was split from, matching the
# <a:rPr lang="it-IT" b="1"/> on the new run in the target markup.
$para2 = $tr.Paragraphs(2, 1)
[void]$para2.InsertAfter("`rTURMA3SI")
